$d = $word.ActiveDocument

# 1) Update the date heading paragraph
$d.Content.Find.Execute("2025-04-11 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-12 Saturday", 2) | Out-Null

# 2) Update each answer cell in the practice table by (row, column) position,
#    since several old values repeat verbatim in different cells.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "82÷8=10, 2"  # was "43÷4=10, 3"
$t.Cell(1, 2).Range.Text = "16÷2=8, 0"  # was "81÷4=20, 1"
$t.Cell(1, 3).Range.Text = "77÷5=15, 2"  # was "80÷4=20, 0"
$t.Cell(1, 4).Range.Text = "75÷8=9, 3"  # was "80÷4=20, 0"
$t.Cell(1, 5).Range.Text = "24÷7=3, 3"  # was "40÷7=5, 5"

$t.Cell(5, 1).Range.Text = "64÷7=9, 1"  # was "81÷4=20, 1"
$t.Cell(5, 2).Range.Text = "94÷9=10, 4"  # was "20÷3=6, 2"
$t.Cell(5, 3).Range.Text = "55÷7=7, 6"  # was "95÷2=47, 1"
$t.Cell(5, 4).Range.Text = "54÷8=6, 6"  # was "59÷7=8, 3"
$t.Cell(5, 5).Range.Text = "83÷3=27, 2"  # was "52÷8=6, 4"

$t.Cell(9, 1).Range.Text = "36÷7=5, 1"  # was "74÷7=10, 4"
$t.Cell(9, 2).Range.Text = "70÷3=23, 1"  # was "11÷6=1, 5"
$t.Cell(9, 3).Range.Text = "27÷6=4, 3"  # was "47÷2=23, 1"
$t.Cell(9, 4).Range.Text = "61÷3=20, 1"  # was "26÷9=2, 8"
$t.Cell(9, 5).Range.Text = "90÷5=18, 0"  # was "32÷7=4, 4"

$t.Cell(13, 1).Range.Text = "21÷5=4, 1"  # was "43÷2=21, 1"
$t.Cell(13, 2).Range.Text = "76÷7=10, 6"  # was "38÷7=5, 3"
$t.Cell(13, 3).Range.Text = "24÷7=3, 3"  # was "68÷8=8, 4"
$t.Cell(13, 4).Range.Text = "30÷3=10, 0"  # was "48÷7=6, 6"
$t.Cell(13, 5).Range.Text = "99÷8=12, 3"  # was "49÷9=5, 4"

$t.Cell(17, 1).Range.Text = "42÷3=14, 0"  # was "41÷3=13, 2"
$t.Cell(17, 2).Range.Text = "17÷6=2, 5"  # was "81÷5=16, 1"
$t.Cell(17, 3).Range.Text = "87÷5=17, 2"  # was "43÷8=5, 3"
$t.Cell(17, 4).Range.Text = "73÷2=36, 1"  # was "11÷3=3, 2"
$t.Cell(17, 5).Range.Text = "33÷4=8, 1"  # was "75÷5=15, 0"

Write-Output "done"
